# Generate Report for Handoff
# Updates the localization-status report: flips the "In Translation" status
# to "Ready for handoff" and refreshes the associated timestamps on all
# three sheets (Overview, zh-cn, de-de), then re-sizes the status columns
# to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: per-language status + the "Latest HO Xliff Generate Date"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 11:03:15"

# --- zh-cn sheet: Status + Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-03 11:03:11"

# --- de-de sheet: Status + Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-03 11:03:15"

# --- Widen the status columns so the longer "Ready for handoff" text fits
# (mirrors Excel's own best-fit column resize after a content update).
$wsOverview.Range("E1").ColumnWidth = 16.333333
$wsOverview.Range("F1").ColumnWidth = 16.333333
$wsZhCn.Range("C1").ColumnWidth = 16.333333
$wsDeDe.Range("C1").ColumnWidth = 16.333333
